$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "ABC company\u2019s management wants to establish an independent and effective internal audit function. They appoint you as a consultant to advise on the authority to whom the head of internal audit should report.Which of the following authorities should you suggest?",
        "ques_type": 2,
        "options": [
            "The audit committee of the board of directors",
            "Chief executive officer",
            "A committee of the company\u2019s senior management",
            "The company\u2019s regulatory body"
        ],
        "score": "The audit committee of the board of directors"
    },
    {
        "title": "You are your company\u2019s internal auditor, and you are considering performing an audit of your company\u2019s sales department. The head of internal audit advises you to perform a detailed risk assessment of the sales function before the start of its audit. Which of the following should you consider the most relevant objective of the activity?",
        "ques_type": 2,
        "options": [
            "To identify audit observations.",
            "To identify resources required for the audit,",
            "To determine risk mitigation strategies.",
            "To determine sample size."
        ],
        "score": "To determine sample size."
    },
    {
        "title": "You are your company\u2019s audit manager, and your team conducts an internal audit of your company\u2019s store outlet on an annual basis. While conducting risk assessment of the outlet, your team observed the following controls:The store incharge prepares bills, receives payments, and reconciles sales with collections each weekend. A security guard stamps each bill when a customer exits the store after shopping.Which of the following control gaps should you consider for mitigation?",
        "ques_type": 2,
        "options": [
            "Inappropriate internal audit frequency",
            "Delayed reconciliation of sales with collections",
            "Manual stamping of bills at the exit gate",
            "Lack of segregation of duties"
        ],
        "score": "Lack of segregation of duties"
    },
    {
        "title": "Your company has a policy to issue a pre-numbered sales invoice and dispatch advice when goods are sold to customers. An outward gate pass is also prepared whenever goods go outside the company\u2019s premises. Customers have the right to reject goods sold if the goods do not meet their requirements. As an internal auditor, you are reviewing a sales transaction to check whether it was correctly recorded as revenue by the accountant.Which of the following documents should you refer to for revenue recognition?",
        "ques_type": 2,
        "options": [
            "Sales invoice",
            "Dispatch advice",
            "Outward gate pass",
            "Customer acknowledgement"
        ],
        "score": "Customer acknowledgement"
    }
]
'@

$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $questionsText
